$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 76-77; existing rows 76-178 shift down to 78-180.
$ws.Rows("76:77").Insert()

# New row 76: Camote, 1a (guarda), Region del Maule, 6/9/2022 (serial 44810)
$ws.Range("A76").Value = 7
$ws.Range("B76").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C76").Value = "Ñuble"
$ws.Range("D76").Value = 44810
$ws.Range("E76").Value = 16
$ws.Range("F76").Value = 100112045
$ws.Range("G76").Value = "Zapallo"
$ws.Range("H76").Value = "Camote"
$ws.Range("I76").Value = "1a (guarda)"
$ws.Range("J76").Value = 300
$ws.Range("K76").Value = 800
$ws.Range("L76").Value = 900
$ws.Range("M76").Value = 850
$ws.Range("N76").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O76").Value = "Región del Maule"
$ws.Range("P76").Value = 850
$ws.Range("Q76").Value = 1
$ws.Range("R76").Value = "Hortaliza"

# New row 77: Camote, 2a (guarda), Region del Maule, 6/9/2022 (serial 44810)
$ws.Range("A77").Value = 7
$ws.Range("B77").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C77").Value = "Ñuble"
$ws.Range("D77").Value = 44810
$ws.Range("E77").Value = 16
$ws.Range("F77").Value = 100112045
$ws.Range("G77").Value = "Zapallo"
$ws.Range("H77").Value = "Camote"
$ws.Range("I77").Value = "2a (guarda)"
$ws.Range("J77").Value = 200
$ws.Range("K77").Value = 700
$ws.Range("L77").Value = 700
$ws.Range("M77").Value = 700
$ws.Range("N77").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O77").Value = "Región del Maule"
$ws.Range("P77").Value = 700
$ws.Range("Q77").Value = 1
$ws.Range("R77").Value = "Hortaliza"
